# Finish grading the "Customer Class" (rows 3-6) and "Product Class" (rows 10-14)
# rubric sections: give full "Points for grading" (column E) equal to the
# "Total Points" already recorded in column D for each graded row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the cursor where the grader last worked (the second section's total).
$ws.Range("E15").Select()
